# Updates Price (D) and Volume(1h) (E) columns for the cryptos list,
# matching the refreshed coinranking.com snapshot.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# row -> @{ D = <new price text, optional>; E = <new volume text> }
$updates = @{
    2 = @{ D='30.051.48'; E='  +7.35%  ' }
    3 = @{ D='1.886.49'; E='  +5.91%  ' }
    4 = @{ D='0.9999'; E='  +0.02%  ' }
    5 = @{ D='249.54'; E='  +2.27%  ' }
    6 = @{ D='0.9998'; E='  +0.04%  ' }
    7 = @{ D='0.4989'; E='  +1.43%  ' }
    8 = @{ D='45.86'; E='  +9.41%  ' }
    9 = @{ D='0.2870'; E='  +7.37%  ' }
    10 = @{ D='0.06592'; E='  +5.31%  ' }
    11 = @{ D='1.881.63'; E='  +5.42%  ' }
    12 = @{ D='17.22'; E='  +5.14%  ' }
    13 = @{ D='0.07215'; E='  +2.60%  ' }
    14 = @{ D='0.6673'; E='  +6.43%  ' }
    15 = @{ D='85.43'; E='  +6.60%  ' }
    16 = @{ D='4.803'; E='  +3.59%  ' }
    17 = @{ D='30.035.57'; E='  +7.42%  ' }
    18 = @{ D='0.9998'; E='  +0.01%  ' }
    19 = @{ D='12.92'; E='  +7.98%  ' }
    20 = @{ D='0.000007526'; E='  +4.12%  ' }
    21 = @{ D='0.9994'; E='  +0.01%  ' }
    22 = @{ D='2.122.16'; E='  +5.52%  ' }
    23 = @{ D='4.781'; E='  +4.52%  ' }
    24 = @{ D='5.539'; E='  +5.85%  ' }
    25 = @{ D='9.044'; E='  +3.86%  ' }
    26 = @{ D='144.50'; E='  +1.87%  ' }
    27 = @{ D='135.48'; E='  +23.67%  ' }
    28 = @{ D='16.74'; E='  +6.21%  ' }
    29 = @{ D='1.976'; E='  +6.17%  ' }
    30 = @{ D='1.395'; E='  +0.47%  ' }
    31 = @{ D='4.183'; E='  -0.31%  ' }
    32 = @{ D='0.08609'; E='  +3.80%  ' }
    33 = @{ D='3.888'; E='  +2.53%  ' }
    34 = @{ E='  +4.04%  ' }
    35 = @{ D='1.145'; E='  +6.86%  ' }
    36 = @{ D='0.6882'; E='  +5.44%  ' }
    37 = @{ E='  +0.07%  ' }
    38 = @{ D='2.710'; E='  +3.79%  ' }
    39 = @{ D='2.312'; E='  +12.79%  ' }
    40 = @{ D='2.755'; E='  +6.66%  ' }
    41 = @{ D='0.9618'; E='  +1.41%  ' }
    42 = @{ D='0.01634'; E='  +5.07%  ' }
    43 = @{ D='6.107'; E='  +2.02%  ' }
    44 = @{ E='  +0.09%  ' }
    45 = @{ D='103.78'; E='  +3.77%  ' }
    46 = @{ D='0.4215'; E='  +5.67%  ' }
    47 = @{ D='7.475'; E='  +3.99%  ' }
    48 = @{ D='0.1256'; E='  +4.43%  ' }
    49 = @{ D='0.05638'; E='  +4.02%  ' }
    50 = @{ D='32.55'; E='  +6.18%  ' }
    51 = @{ D='8.266'; E='  +3.01%  ' }
}

foreach ($row in $updates.Keys) {
    $vals = $updates[$row]
    if ($vals.ContainsKey("D")) {
        # Price column holds plain numeric-looking text (e.g. "0.9999"); force
        # text storage first so Excel does not silently coerce it to a Number,
        # then drop the cell back to its original (unformatted) style.
        $cell = $ws.Range("D$row")
        $cell.NumberFormat = "@"
        $cell.Value = $vals["D"]
        $cell.Style = "Normal"
    }
    $ws.Range("E$row").Value = $vals["E"]
}

Write-Output "Updated cryptos list"
